# Update forecast accuracy numbers on the "Forecasts" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecasts")

# Row 4
$ws.Range("B4").Value = 0.9973886071615455
$ws.Range("E4").Value = 0.01952933962840083

# Row 5
$ws.Range("B5").Value = 68.62732509890299
$ws.Range("E5").Value = 1.24942355082859

# Row 6
$ws.Range("B6").Value = 0.4699252
$ws.Range("C6").Value = 0.5300748
$ws.Range("E6").Value = 0.9981425
$ws.Range("F6").Value = 0.0018575

# Row 7
$ws.Range("B7").Value = 26.5073196
$ws.Range("C7").Value = 26.943518
$ws.Range("E7").Value = 32.5092004
$ws.Range("F7").Value = 16.0649512
